# Publish terminology IG update: bump version, change status, clear the
# "Experimental" flag, update the publication date, and add a Definition
# for the MCS88126 concept (Pt—Bevidsthedsniveau).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 1.0.1
$meta.Range("B3").Value = "1.0.1"

# Status: draft -> active
$meta.Range("B6").Value = "active"

# Experimental: true -> (cleared, no longer set)
$meta.Range("B7").ClearContents()

# Date: 2025-06-28 -> 2025-11-18
# Assigning a date-shaped literal straight to .Value would make Excel
# auto-convert it to a serial date (and pick up a date number format),
# which would change the cell's type/style. Build it as text via TEXT()
# instead, then collapse the formula down to its static text result so
# the stored cell keeps the same "shared string" type and style as
# before.
$dateCell = $meta.Range("B8")
$dateCell.Formula = '=TEXT(DATE(2025,11,18),"yyyy-mm-dd")'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# --- Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Add a Definition for MCS88126 (Pt—Bevidsthedsniveau...) in row 4
$concepts.Range("D4").Value = "0: Habituel, 1: Agiteret eller reagerer kun på tale, 2: Reagerer kun på smerte, 3: Ingen reaktion"
